# Update TPM-derived ligand/receptor expression values for Il15-Il15ra
# (new values recomputed with updated TPM data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.764753333333333
$ws.Range("H2").Value = 14.29426
$ws.Range("I2").Value = 0.2966169987831952
$ws.Range("J2").Value = 0.2966169987831952
$ws.Range("M2").Value = 4.170713666666666
$ws.Range("N2").Value = 12.512141
$ws.Range("O2").Value = 0.1910541549206663
$ws.Range("P2").Value = 0.1910541549206663
$ws.Range("Q2").Value = 19.87242184562889
$ws.Range("R2").Value = 178.85179661066
$ws.Range("S2").Value = 0.05666991003762765
$ws.Range("T2").Value = 0.05666991003762765
$ws.Range("G3").Value = 4.764753333333333
$ws.Range("H3").Value = 14.29426
$ws.Range("I3").Value = 0.2966169987831952
$ws.Range("J3").Value = 0.2966169987831952
$ws.Range("O3").Value = 0.3061812484017604
$ws.Range("P3").Value = 0.3061812484017603
$ws.Range("Q3").Value = 31.84732063004666
$ws.Range("R3").Value = 286.62588567042
$ws.Range("S3").Value = 0.09081856298462214
$ws.Range("T3").Value = 0.09081856298462213
$ws.Range("G4").Value = 4.764753333333333
$ws.Range("H4").Value = 14.29426
$ws.Range("I4").Value = 0.2966169987831952
$ws.Range("J4").Value = 0.2966169987831952
$ws.Range("M4").Value = 10.975355
$ws.Range("N4").Value = 32.926065
$ws.Range("O4").Value = 0.5027645966775733
$ws.Range("P4").Value = 0.5027645966775732
$ws.Range("Q4").Value = 52.29485932076667
$ws.Range("R4").Value = 470.6537338869
$ws.Range("S4").Value = 0.1491285257609454
$ws.Range("T4").Value = 0.1491285257609454
$ws.Range("I5").Value = 0.5337607564504776
$ws.Range("J5").Value = 0.5337607564504775
$ws.Range("M5").Value = 4.170713666666666
$ws.Range("N5").Value = 12.512141
$ws.Range("O5").Value = 0.1910541549206663
$ws.Range("P5").Value = 0.1910541549206663
$ws.Range("Q5").Value = 35.76032041433633
$ws.Range("R5").Value = 321.8428837290269
$ws.Range("S5").Value = 0.1019772102534616
$ws.Range("T5").Value = 0.1019772102534615
$ws.Range("I6").Value = 0.5337607564504776
$ws.Range("J6").Value = 0.5337607564504775
$ws.Range("O6").Value = 0.3061812484017604
$ws.Range("P6").Value = 0.3061812484017603
$ws.Range("S6").Value = 0.1634275347578752
$ws.Range("T6").Value = 0.1634275347578751
$ws.Range("I7").Value = 0.5337607564504776
$ws.Range("J7").Value = 0.5337607564504775
$ws.Range("M7").Value = 10.975355
$ws.Range("N7").Value = 32.926065
$ws.Range("O7").Value = 0.5027645966775733
$ws.Range("P7").Value = 0.5027645966775732
$ws.Range("Q7").Value = 94.104329097895
$ws.Range("R7").Value = 846.938961881055
$ws.Range("S7").Value = 0.2683560114391408
$ws.Range("T7").Value = 0.2683560114391407
$ws.Range("G8").Value = 2.724753333333334
$ws.Range("H8").Value = 8.17426
$ws.Range("I8").Value = 0.1696222447663273
$ws.Range("J8").Value = 0.1696222447663273
$ws.Range("M8").Value = 4.170713666666666
$ws.Range("N8").Value = 12.512141
$ws.Range("O8").Value = 0.1910541549206663
$ws.Range("P8").Value = 0.1910541549206663
$ws.Range("Q8").Value = 11.36416596562889
$ws.Range("R8").Value = 102.27749369066
$ws.Range("S8").Value = 0.03240703462957707
$ws.Range("T8").Value = 0.03240703462957706
$ws.Range("G9").Value = 2.724753333333334
$ws.Range("H9").Value = 8.17426
$ws.Range("I9").Value = 0.1696222447663273
$ws.Range("J9").Value = 0.1696222447663273
$ws.Range("O9").Value = 0.3061812484017604
$ws.Range("P9").Value = 0.3061812484017603
$ws.Range("Q9").Value = 18.21208507004667
$ws.Range("R9").Value = 163.90876563042
$ws.Range("S9").Value = 0.05193515065926305
$ws.Range("T9").Value = 0.05193515065926303
$ws.Range("G10").Value = 2.724753333333334
$ws.Range("H10").Value = 8.17426
$ws.Range("I10").Value = 0.1696222447663273
$ws.Range("J10").Value = 0.1696222447663273
$ws.Range("M10").Value = 10.975355
$ws.Range("N10").Value = 32.926065
$ws.Range("O10").Value = 0.5027645966775733
$ws.Range("P10").Value = 0.5027645966775732
$ws.Range("Q10").Value = 29.90513512076667
$ws.Range("R10").Value = 269.1462160869
$ws.Range("S10").Value = 0.08528005947748717
$ws.Range("T10").Value = 0.08528005947748712
